$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 18 (meeting #3, second occurrence) with the new meeting
# date/time, place, and attendance marks for the team - matching the
# pattern already used for rows 10-17.
$ws.Range("B18").Value = "10/3 /1:00"
$ws.Range("C18").Value = "Google Hangout"
$ws.Range("D18").Value = "A"
$ws.Range("E18").Value = "A"
$ws.Range("F18").Value = "T"
$ws.Range("G18").Value = "T"
$ws.Range("H18").Value = "A"
$ws.Range("I18").Value = "A"

# Row 18 previously used the "last row" border style (like C19:C21, etc.);
# bring it in line with the rest of the data rows (10-17) by copying the
# cell formatting from the row above (row 17), which already has that look.
$ws.Range("C17").Copy()
$ws.Range("C18").PasteSpecial(-4122)

# Update the current selection/view to match where the user ended up after
# editing - cell B19, with no special top-left scroll anchor.
$ws.Range("B19").Select()
